# Weekly refresh of the "Rabanito" sheet: the underlying data rows (2-43)
# were re-sorted / re-shuffled when the new week's data was merged in.
# The net effect (confirmed against the canonical XML diff) is a pure
# permutation of whole data rows A:R for rows 2..43 - i.e. every "after"
# row is an exact copy of some "before" row, just relocated.
#
# Build the permutation as: new row <N> <- old row <mapping[N]>.
$map = @{
    2=2; 3=31; 4=32; 5=22; 6=6; 7=15; 8=3; 9=19; 10=18; 11=10;
    12=34; 13=20; 14=30; 15=16; 16=8; 17=40; 18=23; 19=24; 20=7; 21=36;
    22=13; 23=4; 24=41; 25=27; 26=12; 27=11; 28=5; 29=37; 30=42; 31=9;
    32=14; 33=25; 34=26; 35=43; 36=38; 37=21; 38=35; 39=39; 40=17; 41=33;
    42=28; 43=29
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every data row (A:R) BEFORE any writes, so overwriting a row
# never clobbers data we still need to read for a later target row.
$snapshot = @{}
for ($r = 2; $r -le 43; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":R" + $r).Value()
}

# Now write each target row from its mapped source snapshot.
foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $ws.Range("A" + $newRow + ":R" + $newRow).Value = $snapshot[$oldRow]
}
